# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Updates the "K" column (column G) values for the pitching log on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new K (strikeout) value
$kValues = @{
    2  = 0
    3  = 1
    4  = 0
    5  = 3
    6  = 2
    7  = 2
    8  = 1
    9  = 3
    10 = 6
    11 = 5
    12 = 4
    13 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
